# Included missing group in DIF analyses for migration in data example 1.
#
# The "mig" (migration) DIF analysis previously dropped one migration
# group, so the "estimates" sheet only reported a single "mig 0-1"
# contrast and the "gof" sheet's mig rows were fit on a smaller sample
# (N=814) than the "sex" rows (N=1000). Re-running with all groups
# included adds two more pairwise contrasts ("mig 1-3", "mig 2-3") and
# brings the mig rows up to the full N=1000 sample, which shifts the
# goodness-of-fit statistics and several of the previously reported
# estimates.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "gof"
$ws2 = $wb.Worksheets.Item(2)   # "estimates"

# ----------------------------------------------------------------------
# Sheet "gof": the "mig" rows (4 & 5) now use the full N=1000 sample,
# changing Deviance / Number.of.parameters / AIC / BIC.
# ----------------------------------------------------------------------
$ws1.Range("C4:G5").Value = @(
    @(1000, 13597, 18, 13633, 13721),
    @(1000, 13570, 46, 13662, 13888)
)

# ----------------------------------------------------------------------
# Sheet "estimates": header C1 ("mig 0-1") becomes "mig 1-2", and two
# new columns are added for the previously-missing contrasts.
# ----------------------------------------------------------------------
$ws2.Range("C1").Value = "mig 1-2"
$ws2.Range("D1").Value = "mig 1-3"
$ws2.Range("E1").Value = "mig 2-3"

$data = @(
    @(" 0.017 ( 0.015)", "-0.519 (-0.444)", "-0.536 (-0.459)"),  # grk10001_c
    @("-0.023 (-0.020)", " 0.226 ( 0.193)", " 0.250 ( 0.214)"),  # grk10002_c
    @("-0.064 (-0.055)", "-0.137 (-0.117)", "-0.073 (-0.062)"),  # grk10003_c
    @(" 0.253 ( 0.217)", " 0.244 ( 0.209)", "-0.009 (-0.008)"),  # grk10004_c
    @("-0.259 (-0.222)", " 0.345 ( 0.295)", " 0.603 ( 0.516)"),  # grk10005_c
    @(" 0.094 ( 0.080)", " 0.359 ( 0.307)", " 0.265 ( 0.227)"),  # grk10006_c
    @(" 0.078 ( 0.067)", " 0.164 ( 0.140)", " 0.086 ( 0.074)"),  # grk10007_c
    @(" 0.096 ( 0.082)", "-0.026 (-0.022)", "-0.122 (-0.104)"),  # grk10008_c
    @("-0.071 (-0.061)", " 0.110 ( 0.094)", " 0.181 ( 0.155)"),  # grk10009_c
    @(" 0.054 ( 0.046)", "-0.112 (-0.096)", "-0.166 (-0.142)"),  # grk10010_c
    @(" 0.181 ( 0.155)", "-0.019 (-0.016)", "-0.199 (-0.170)"),  # grk10011_c
    @(" 0.141 ( 0.121)", "-0.094 (-0.080)", "-0.235 (-0.201)"),  # grk10012_c
    @(" 0.090 ( 0.077)", "-0.115 (-0.098)", "-0.205 (-0.175)"),  # grk10013_c
    @("-0.131 (-0.112)", " 0.026 ( 0.022)", " 0.157 ( 0.134)"),  # grk10014_c
    @(" 0.456 ( 0.390)", " 0.452 ( 0.387)", "-0.004 (-0.003)"),  # grk10015_c
    @("0.756 (0.647)",   "0.538 (0.461)",   "-0.218 (-0.187)"),  # Main effect (DIF model)
    @("0.730 (0.627)",   "0.497 (0.427)",   "-0.233 (-0.200)")   # Main effect (Main effect model)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws2.Range("C$row").Value = $data[$i][0]
    $ws2.Range("D$row").Value = $data[$i][1]
    $ws2.Range("E$row").Value = $data[$i][2]
}
